$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 13-21: plain data rows (same look as rows 9-11) ---
$ws.Range("A13").Value = 'SCRIPT/T01P01A/um1105.ssb'
$ws.Range("B13").Value = 427
$ws.Range("C13").Value = ' Oh, [hero] and\n[partner]!'
$ws.Range("D13").Value = ' О, [hero] и [partner]!'
$ws.Range("E13").Value = ' Ï, [hero] é [partner]!'

$ws.Range("A14").Value = 'SCRIPT/T01P01A/um1108.ssb'
$ws.Range("B14").Value = 430
$ws.Range("C14").Value = ' I\''ve come into an interesting bit\nof news I can\''t wait to share with you.'
$ws.Range("D14").Value = ' Я тут узнал кое-что очень\nинтересное и хочу вам это рассказать.'
$ws.Range("E14").Value = ' Ÿ óôó ôèîàì ëïå-œóï ïœåîû\néîóåñåòîïå é öïœô âàí üóï ñàòòëàèàóû.'

$ws.Range("A15").Value = 'SCRIPT/T01P01A/um1111.ssb'
$ws.Range("B15").Value = 433
$ws.Range("C15").Value = ' Have you heard about Eggs,\n[hero]?'
$ws.Range("D15").Value = ' Ты знаешь о Яйцах,\n[hero]?'
$ws.Range("E15").Value = ' Óú èîàåšû ï Ÿêøàö,\n[hero]?'

$ws.Range("B16").Value = 436
$ws.Range("C16").Value = ' Eggs are quite rare, I hear.[K] And\nI hear you may get them as rewards for jobs.'
$ws.Range("D16").Value = ' Яйца довольно редко попадаются.\nГоворят, что их можно получить как\nнаграду за задание.'
$ws.Range("E16").Value = ' Ÿêøà äïâïìûîï ñåäëï ðïðàäàýóòÿ.\nÃïâïñÿó, œóï éö íïçîï ðïìôœéóû ëàë\nîàãñàäô èà èàäàîéå.'

$ws.Range("B17").Value = 439
$ws.Range("C17").Value = ' But as precious as Eggs can be,\nthey can be a burden to look after.'
$ws.Range("D17").Value = ' Яйца конечно ценны, но за ними\nочень трудно следить.'
$ws.Range("E17").Value = ' Ÿêøà ëïîåœîï øåîîú, îï èà îéíé\nïœåîû óñôäîï òìåäéóû.'

$ws.Range("B18").Value = 442
$ws.Range("C18").Value = ' But...[K]worry no more!'
$ws.Range("D18").Value = ' Но...[K] Нет волненью!'
$ws.Range("E18").Value = ' Îï...[K] Îåó âïìîåîûý!'

$ws.Range("B19").Value = 445
$ws.Range("C19").Value = ' There\''s a new shop in town that\nwill take care of your Eggs!'
$ws.Range("D19").Value = ' В городе появилось новое место,\nгде могут позаботиться о твоих Яйцах!'
$ws.Range("E19").Value = ' Â ãïñïäå ðïÿâéìïòû îïâïå íåòóï,\nãäå íïãôó ðïèàáïóéóûòÿ ï óâïéö Ÿêøàö!'

$ws.Range("B20").Value = 461
$ws.Range("C20").Value = ' See that shop in the lower left?'
$ws.Range("D20").Value = ' Видите его слева внизу?'
$ws.Range("E20").Value = ' Âéäéóå åãï òìåâà âîéèô?'

$ws.Range("B21").Value = 476
$ws.Range("C21").Value = ' A Pokémon named [CS:N]Chansey[CR] runs\nthe shop.'
$ws.Range("D21").Value = ' Этим место управляет Покемон\nпо имени [CS:N]Ченси[CR].'
$ws.Range("E21").Value = ' Üóéí íåòóï ôðñàâìÿåó Ðïëåíïî\nðï éíåîé [CS:N]Œåîòé[CR].'

# --- Row 22: new data row (values first, section-ending format applied below with row 12) ---
$ws.Range("B22").Value = 479
$ws.Range("C22").Value = ' If you ever get an Egg, [CS:N]Chansey[CR]\nwill take care of it for you, [hero].'
$ws.Range("D22").Value = ' Если ты получишь Яйцо, [CS:N]Ченси[CR]\nбудет о нём заботиться, [hero].'
$ws.Range("E22").Value = ' Åòìé óú ðïìôœéšû Ÿêøï, [CS:N]Œåîòé[CR]\náôäåó ï îæí èàáïóéóûòÿ, [hero].'

# --- Section-ending (bottom border) look: row 12 becomes mid-table, row 22 becomes the new last row ---
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A12:E12").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A22:E22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row heights (Excel auto-calculates these on real wrap; set explicitly to match) ---
$ws.Rows("13:13").RowHeight = 43.2
$ws.Rows("14:14").RowHeight = 43.2
$ws.Rows("15:15").RowHeight = 43.2
$ws.Rows("16:16").RowHeight = 31.8
$ws.Rows("17:17").RowHeight = 21.6
$ws.Rows("19:19").RowHeight = 21.6
$ws.Rows("21:21").RowHeight = 21.6
$ws.Rows("22:22").RowHeight = 31.8

# --- View state: scrolled down, new active selection ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D15").Select() | Out-Null
